$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Record"
$ws.Range("B10").Value = "Balanço Geral"
$ws.Range("C10").Value = "Limpeza Pública"
$ws.Range("D10").Value = "2025-04-02T13:27"
$ws.Range("E10").Value = "Neutro"
$ws.Range("F10").Value = "Tô na Bronca. Árvore obstrui iluminação de poste e moradores solicitam poda. É na Rua José Caldas de Carvalho, no Jóquei 2. *com nota* Equipe da secretaria de Serviços Públicos vai avaliar se é responsabilidade da prefeitura. "

$ws.Range("A11").Value = "Record"
$ws.Range("B11").Value = "Balanço Geral"
$ws.Range("C11").Value = "Infraestrutura"
$ws.Range("D11").Value = "2025-04-02T13:29"
$ws.Range("E11").Value = "Neutro"
$ws.Range("F11").Value = "Tô na bronca. Pedestres cobram melhorias na passarela sobre a RJ 216 no Parque Imperial. Vídeo exibido por um morador. Muito mato, sem acesso à rampa. Prefeitura tem mover ação junto ao Governo Estadual para solucionar o problema. DER enviou nota. "

$ws.Range("A12").Value = "Record"
$ws.Range("B12").Value = "Balanço Geral"
$ws.Range("C12").Value = "Social"
$ws.Range("D12").Value = "2025-04-02T13:45"
$ws.Range("E12").Value = "Neutro"
$ws.Range("F12").Value = "Dia do Autismo. Data reforça a conscientização mundial sobre o transtorno do espectro autista. Entrevista com mãe e com presidente da Apape, Naira Peçanha. *matéria*"

$ws.Range("A13").Value = "Record"
$ws.Range("B13").Value = "Balanço Geral"
$ws.Range("C13").Value = "Governo"
$ws.Range("D13").Value = "2025-04-02T11:51"
$ws.Range("E13").Value = "Positivo"
$ws.Range("F13").Value = "Projeto de reforma administrativa começa a tramitar na sessão de hoje na Câmara. *nota coberta*"
